$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (the quarterly financials gain
# two new reporting periods: 2018-12-31 and 2018-09-30). Existing columns
# D:K shift right to F:M.
$ws.Columns("D:E").Insert()

# The newly inserted D:E columns land with default (unstyled) formatting.
# Copy the number formatting from column F (the old column D, now shifted
# right, which already carries the correct style) into the new columns so
# they match the rest of the table. This is done per contiguous data block
# so that blank separator rows (5:6, 36:37, 79) are left untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("F7:F35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new quarterly columns with their reported figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 233200
$ws.Range("E8").Value = 176900
$ws.Range("D9").Value = 147500
$ws.Range("E9").Value = 108900
$ws.Range("D10").Value = 85700
$ws.Range("E10").Value = 68000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 4500
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 10200
$ws.Range("E15").Value = 10000
$ws.Range("D17").Value = 207300
$ws.Range("E17").Value = 160100
$ws.Range("D18").Value = 25900
$ws.Range("E18").Value = 16800
$ws.Range("D20").Value = -1500
$ws.Range("E20").Value = -1400
$ws.Range("D21").Value = 34600
$ws.Range("E21").Value = 25300
$ws.Range("D22").Value = 4500
$ws.Range("E22").Value = 5800
$ws.Range("D23").Value = 19900
$ws.Range("E23").Value = 9600
$ws.Range("D24").Value = 2800
$ws.Range("E24").Value = 1500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 17100
$ws.Range("E26").Value = 8100
$ws.Range("D27").Value = 6000
$ws.Range("E27").Value = 2000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 1500
$ws.Range("E32").Value = 1400
$ws.Range("D33").Value = 6000
$ws.Range("E33").Value = 2000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 6000
$ws.Range("E35").Value = 2000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 13500
$ws.Range("E41").Value = 10600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 148600
$ws.Range("E43").Value = 127000
$ws.Range("D44").Value = 86600
$ws.Range("E44").Value = 81200
$ws.Range("D45").Value = 11900
$ws.Range("E45").Value = 25300
$ws.Range("D46").Value = 260600
$ws.Range("E46").Value = 244100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 44300
$ws.Range("E48").Value = 44500
$ws.Range("D49").Value = 346500
$ws.Range("E49").Value = 350500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 11600
$ws.Range("E52").Value = 9700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 663000
$ws.Range("E54").Value = 648800
$ws.Range("D57").Value = 36100
$ws.Range("E57").Value = 35000
$ws.Range("D58").Value = 30600
$ws.Range("E58").Value = 62900
$ws.Range("D59").Value = 71100
$ws.Range("E59").Value = 67400
$ws.Range("D60").Value = 137900
$ws.Range("E60").Value = 165300
$ws.Range("D61").Value = 216700
$ws.Range("E61").Value = 196200
$ws.Range("D62").Value = 12100
$ws.Range("E62").Value = 10900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 506400
$ws.Range("E66").Value = 501800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 10400
$ws.Range("E72").Value = 4400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 156600
$ws.Range("E76").Value = 147000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 6000
$ws.Range("E81").Value = 2000
$ws.Range("D83").Value = 10200
$ws.Range("E83").Value = 10000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 28600
$ws.Range("E89").Value = -700
$ws.Range("D91").Value = -6000
$ws.Range("E91").Value = -6800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -6000
$ws.Range("E94").Value = -6800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -18700
$ws.Range("E100").Value = 8700
$ws.Range("D101").Value = -1000
$ws.Range("E101").Value = -1500
$ws.Range("D102").Value = 2900
$ws.Range("E102").Value = -300
